$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 11.11392488232964
$ws.Range("C2").Value = 8.486287452005675
$ws.Range("D2").Value = 3.893111034493049
$ws.Range("F2").Value = 20.54983768875346
$ws.Range("G2").Value = 3.598744595220072
$ws.Range("M2").Value = 19.96815084659818
$ws.Range("O2").Value = 18.22946327563676

$ws.Range("B3").Value = 10.54320204495148
$ws.Range("C3").Value = 8.194613980783108
$ws.Range("D3").Value = 3.843837355465287
$ws.Range("F3").Value = 20.51677351315613
$ws.Range("G3").Value = 3.600876652675796
$ws.Range("M3").Value = 19.36275468584799
$ws.Range("O3").Value = 18.27408084785947

$ws.Range("B4").Value = 10.17664796610747
$ws.Range("C4").Value = 8.009102710566873
$ws.Range("D4").Value = 3.812968854709414
$ws.Range("F4").Value = 20.5045364596227
$ws.Range("G4").Value = 3.602254047048987
$ws.Range("M4").Value = 18.98949138819381
$ws.Range("O4").Value = 18.30770238345032

$ws.Range("B5").Value = 10.02337290749713
$ws.Range("C5").Value = 7.931975168187614
$ws.Range("D5").Value = 3.800243612340682
$ws.Range("F5").Value = 20.50157707320447
$ws.Range("G5").Value = 3.602832581133589
$ws.Range("M5").Value = 18.83725019432145
$ws.Range("O5").Value = 18.32296017933113

$ws.Range("B6").Value = 9.997691131159248
$ws.Range("C6").Value = 7.919078435213995
$ws.Range("D6").Value = 3.798122021568361
$ws.Range("F6").Value = 20.50120804615518
$ws.Range("G6").Value = 3.602929688957185
$ws.Range("M6").Value = 18.81197032764604
$ws.Range("O6").Value = 18.32558747580883

$ws.Range("B7").Value = 10.17459641703237
$ws.Range("C7").Value = 8.008068618551398
$ws.Range("D7").Value = 3.812797817890311
$ws.Range("F7").Value = 20.50448834242157
$ws.Range("G7").Value = 3.602261779499376
$ws.Range("M7").Value = 18.98743840280071
$ws.Range("O7").Value = 18.30790186438637

$ws.Range("B8").Value = 10.92057120992136
$ws.Range("C8").Value = 8.3870995599955
$ws.Range("D8").Value = 3.876253925600098
$ws.Range("F8").Value = 20.5367638261656
$ws.Range("G8").Value = 3.599465588349669
$ws.Range("M8").Value = 19.7598807839736
$ws.Range("O8").Value = 18.24355004268488

$ws.Range("B9").Value = 12.24997259504609
$ws.Range("C9").Value = 9.076138378953701
$ws.Range("D9").Value = 3.995460537731337
$ws.Range("F9").Value = 20.66394696462158
$ws.Range("G9").Value = 3.594521495575549
$ws.Range("M9").Value = 21.25184238515204
$ws.Range("O9").Value = 18.16712312945393

$ws.Range("B10").Value = 13.13949325327647
$ws.Range("C10").Value = 9.545500915666079
$ws.Range("D10").Value = 4.079392928661238
$ws.Range("F10").Value = 20.79603049118422
$ws.Range("G10").Value = 3.591213983218111
$ws.Range("M10").Value = 22.32115074527741
$ws.Range("O10").Value = 18.14176204469386

$ws.Range("B11").Value = 13.52434934546593
$ws.Range("C11").Value = 9.750366173547343
$ws.Range("D11").Value = 4.116691390169731
$ws.Range("F11").Value = 20.86438781623001
$ws.Range("G11").Value = 3.589779046248619
$ws.Range("M11").Value = 22.79939139864959
$ws.Range("O11").Value = 18.13698988187796

$ws.Range("B12").Value = 13.66718332017413
$ws.Range("C12").Value = 9.826653897370463
$ws.Range("D12").Value = 4.130681025740805
$ws.Range("F12").Value = 20.8914480952364
$ws.Range("G12").Value = 3.589245628626381
$ws.Range("M12").Value = 22.97913566985564
$ws.Range("O12").Value = 18.13616065199404

$ws.Range("B13").Value = 13.63655134854302
$ws.Range("C13").Value = 9.810282033072982
$ws.Range("D13").Value = 4.127674205222469
$ws.Range("F13").Value = 20.88556818850657
$ws.Range("G13").Value = 3.589360067490655
$ws.Range("M13").Value = 22.94048742181025
$ws.Range("O13").Value = 18.1362956780759

$ws.Range("B14").Value = 13.53615884103217
$ws.Range("C14").Value = 9.756668514630642
$ws.Range("D14").Value = 4.117845063193765
$ws.Range("F14").Value = 20.86659062892209
$ws.Range("G14").Value = 3.589734962329932
$ws.Range("M14").Value = 22.81420717639736
$ws.Range("O14").Value = 18.13690203632164

$ws.Range("B15").Value = 13.47428602893259
$ws.Range("C15").Value = 9.723659376210261
$ws.Range("D15").Value = 4.111806697965482
$ws.Range("F15").Value = 20.85511884718386
$ws.Range("G15").Value = 3.589965891910048
$ws.Range("M15").Value = 22.73667545941336
$ws.Range("O15").Value = 18.1374009339629

$ws.Range("B16").Value = 13.11393863755591
$ws.Range("C16").Value = 9.531934318480623
$ws.Range("D16").Value = 4.07693697093081
$ws.Range("F16").Value = 20.79172843707962
$ws.Range("G16").Value = 3.591309156795211
$ws.Range("M16").Value = 22.28971619761518
$ws.Range("O16").Value = 18.14221049953908

$ws.Range("B17").Value = 12.88776404802977
$ws.Range("C17").Value = 9.412066709619243
$ws.Range("D17").Value = 4.055313955746486
$ws.Range("F17").Value = 20.75494949416591
$ws.Range("G17").Value = 3.592151009648152
$ws.Range("M17").Value = 22.01329117875635
$ws.Range("O17").Value = 18.14689755920955

$ws.Range("B18").Value = 12.75581413125818
$ws.Range("C18").Value = 9.342310207106229
$ws.Range("D18").Value = 4.042794307068281
$ws.Range("F18").Value = 20.73457484829633
$ws.Range("G18").Value = 3.592641781873441
$ws.Range("M18").Value = 21.85353708350743
$ws.Range("O18").Value = 18.15022991091718

$ws.Range("B19").Value = 12.71082057267353
$ws.Range("C19").Value = 9.318553937668968
$ws.Range("D19").Value = 4.038541406499045
$ws.Range("F19").Value = 20.72781064946111
$ws.Range("G19").Value = 3.592809077443101
$ws.Range("M19").Value = 21.79932194582982
$ws.Range("O19").Value = 18.15146733244268

$ws.Range("B20").Value = 12.9120336979243
$ws.Range("C20").Value = 9.424911197565713
$ws.Range("D20").Value = 4.057624373428274
$ws.Range("F20").Value = 20.75878408472758
$ws.Range("G20").Value = 3.592060714371364
$ws.Range("M20").Value = 22.04279726873576
$ws.Range("O20").Value = 18.1463327004984

$ws.Range("B21").Value = 13.565725727592
$ws.Range("C21").Value = 9.772451463211635
$ws.Range("D21").Value = 4.120735828794088
$ws.Range("F21").Value = 20.87213304251047
$ws.Range("G21").Value = 3.589624576712605
$ws.Range("M21").Value = 22.85133684439636
$ws.Range("O21").Value = 18.13669735845511

$ws.Range("B22").Value = 13.97601356211137
$ws.Range("C22").Value = 9.99205083793194
$ws.Range("D22").Value = 4.161195673280252
$ws.Range("F22").Value = 20.95305247043708
$ws.Range("G22").Value = 3.588090459396454
$ws.Range("M22").Value = 23.37179604730623
$ws.Range("O22").Value = 18.13610168082096

$ws.Range("B23").Value = 13.7586009115933
$ws.Range("C23").Value = 9.875550079765434
$ws.Range("D23").Value = 4.139675924960154
$ws.Range("F23").Value = 20.90924393539333
$ws.Range("G23").Value = 3.588903954562121
$ws.Range("M23").Value = 23.09480016075598
$ws.Range("O23").Value = 18.13589645867537

$ws.Range("B24").Value = 12.90106736245833
$ws.Range("C24").Value = 9.419106826790223
$ws.Range("D24").Value = 4.056580108020767
$ws.Range("F24").Value = 20.75704806669496
$ws.Range("G24").Value = 3.592101515744061
$ws.Range("M24").Value = 22.02946014853265
$ws.Range("O24").Value = 18.14658608723218

$ws.Range("B25").Value = 11.90529723621311
$ws.Range("C25").Value = 8.895965453328493
$ws.Range("D25").Value = 3.963820444811897
$ws.Range("F25").Value = 20.62272170234496
$ws.Range("G25").Value = 3.59580167082497
$ws.Range("M25").Value = 20.85206972083139
$ws.Range("O25").Value = 18.18242108371829

